$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns A (Company Name), B (Company Number), H (Category)
# for rows 2-8, reflecting the re-ordering described in the diff.
$values = @{
    2 = @("T GILPIN PHYSIO CONSULTANCY LTD", "16460503", "LP")
    3 = @("SAMVIV PARTNERS LTD", "16460672", "Partners")
    4 = @("4D CAPITAL PROPCO (44) LIMITED", "16461269", "Capital")
    5 = @("DGPI LTD", "SC849118", "GP")
    6 = @("DAVIDSON CAPITAL HOLDINGS LTD", "SC849117", "Capital")
    7 = @("AFROSCOT VENTURES LTD", "16462878", "Ventures")
    8 = @("ST GEORGE CAPITAL (LAND) LIMITED", "16462880", "Capital")
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Range("A$row").Value = $v[0]
    # Company numbers are stored as text (some contain letters, e.g. Scottish
    # numbers like "SC849118"). Force text format only when the new value is
    # purely numeric, otherwise it would be auto-converted to a number.
    if ($v[1] -match '^\d+$') {
        $ws.Range("B$row").NumberFormat = "@"
    }
    $ws.Range("B$row").Value = $v[1]
    $ws.Range("H$row").Value = $v[2]
}

$wb.Save()
